$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "42.378.60" using
# "." as a thousands separator) - force text format so Excel COM does not
# silently coerce them into Double values / normalize their formatting.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.378.60'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.285.94'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '157.04'
$ws.Range('E5').Value = '  +15,591.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '307.07'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '95.53'
$ws.Range('E7').Value = '  +4.60%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '35.79'
$ws.Range('E11').Value = '  +10.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0802'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.73'
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.641.70'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.45'
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.292.40'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('E18').Value = '  +4.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.291.03'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.66'
$ws.Range('E20').Value = '  +4.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0916'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.01'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '242.82'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.60'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.09'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.08'
$ws.Range('E29').Value = '  +5.29%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('E31').Value = '  -8.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.87'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.33'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0753'
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('E36').Value = '  +2.62%  '
$ws.Range('E37').Value = '  +5.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.24'
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('E42').Value = '  +7.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.012.06'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.36'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  +10.95%  '
$ws.Range('E46').Value = '  +2.40%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.99'
$ws.Range('E48').Value = '  +4.59%  '
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.41'
$ws.Range('E50').Value = '  +3.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.92'
$ws.Range('E51').Value = '  -0.56%  '
